# Subjective Evaluation (Metadata) - fill in reviewer (RB) scores for the
# four tracks being evaluated. Each row corresponds to one ranked result
# (B3:B22 = 1..20); columns C/H/M/R hold the "RB" rating for each of the
# four tracks (the SS/TV columns - D/E, I/J, N/O, S/T - are left blank,
# same as before). The Mean/Std./Global-mean/Global-std. columns are all
# formulas already in the sheet, so they recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, C, H, M, R
$data = @(
    @(3,  3,3,3,5),
    @(4,  2,3,2,5),
    @(5,  5,4,4,5),
    @(6,  3,2,5,5),
    @(7,  3,2,5,5),
    @(8,  1,1,4,1),
    @(9,  3,3,4,5),
    @(10, 2,1,3,1),
    @(11, 1,5,1,1),
    @(12, 2,2,4,1),
    @(13, 1,2,4,1),
    @(14, 1,4,1,1),
    @(15, 2,4,3,1),
    @(16, 1,4,2,4),
    @(17, 3,3,2,5),
    @(18, 2,3,5,5),
    @(19, 1,3,1,2),
    @(20, 1,2,2,2),
    @(21, 2,5,4,4),
    @(22, 3,4,4,4)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Range("C$r").Value = $entry[1]
    $ws.Range("H$r").Value = $entry[2]
    $ws.Range("M$r").Value = $entry[3]
    $ws.Range("R$r").Value = $entry[4]
}

# Recalculate so the Mean/Std/threshold formulas pick up the new values.
$excel.Calculate()

# Leave the sheet scrolled/selected where the author left off.
$ws.Range("T25").Select()
